# Rewrites the body of the document to match the new db_handler.py listing.
# The whole-document Range.Text assignment (joined with CR paragraph marks)
# replaces every paragraph in one shot, which is far more reliable than trying
# to chase each Find/Replace hunk through a document whose paragraph count and
# ordering changed substantially.
$lines = @(
    "import sqlite3",
    "",
    "def connect():",
    "    conn = sqlite3.connect(`"hostel.db`")",
    "    return conn",
    "",
    "def create_tables():",
    "    conn = connect()",
    "    cursor = conn.cursor()",
    "    # Students table",
    "    cursor.execute(`"`"`"",
    "        CREATE TABLE IF NOT EXISTS students (",
    "            id TEXT PRIMARY KEY,",
    "            name TEXT,",
    "            age INTEGER,",
    "            room TEXT",
    "        )",
    "    `"`"`")",
    "    # Rooms table",
    "    cursor.execute(`"`"`"",
    "        CREATE TABLE IF NOT EXISTS rooms (",
    "            room_no TEXT PRIMARY KEY,",
    "            capacity INTEGER,",
    "            occupants INTEGER DEFAULT 0",
    "        )",
    "    `"`"`")",
    "    # Payments table",
    "    cursor.execute(`"`"`"",
    "        CREATE TABLE IF NOT EXISTS payments (",
    "            payment_id INTEGER PRIMARY KEY AUTOINCREMENT,",
    "            student_id TEXT,",
    "            amount REAL,",
    "            date TEXT",
    "        )",
    "    `"`"`")",
    "    # Insert rooms 1-100 if not exists",
    "    for i in range(1, 101):",
    "        cursor.execute(`"INSERT OR IGNORE INTO rooms (room_no, capacity, occupants) VALUES (?, 2, 0)`", (str(i),))",
    "    conn.commit()",
    "    conn.close()",
    "",
    "def assign_room():",
    "    conn = connect()",
    "    cursor = conn.cursor()",
    "    cursor.execute(`"SELECT room_no FROM rooms WHERE occupants < capacity ORDER BY room_no ASC LIMIT 1`")",
    "    room = cursor.fetchone()",
    "    conn.close()",
    "    if room:",
    "        return room[0]",
    "    else:",
    "        return None",
    "",
    "def add_student(id, name, age):",
    "    conn = connect()",
    "    cursor = conn.cursor()",
    "    room = assign_room()",
    "    if not room:",
    "        raise Exception(`"No available rooms!`")",
    "    try:",
    "        cursor.execute(`"INSERT INTO students (id, name, age, room) VALUES (?, ?, ?, ?)`", (id, name, age, room))",
    "        cursor.execute(`"UPDATE rooms SET occupants = occupants + 1 WHERE room_no = ?`", (room,))",
    "        conn.commit()",
    "    except sqlite3.IntegrityError:",
    "        raise Exception(`"Student ID already exists!`")",
    "    conn.close()",
    "",
    "def delete_student(student_id):",
    "    conn = connect()",
    "    cursor = conn.cursor()",
    "    cursor.execute(`"SELECT room FROM students WHERE id = ?`", (student_id,))",
    "    room = cursor.fetchone()",
    "    if room:",
    "        cursor.execute(`"DELETE FROM students WHERE id = ?`", (student_id,))",
    "        cursor.execute(`"UPDATE rooms SET occupants = occupants - 1 WHERE room_no = ?`", (room[0],))",
    "        conn.commit()",
    "    conn.close()",
    "",
    "def view_students():",
    "    conn = connect()",
    "    cursor = conn.cursor()",
    "    cursor.execute(`"SELECT * FROM students`")",
    "    rows = cursor.fetchall()",
    "    conn.close()",
    "    return rows",
    "",
    "def search_student(keyword):",
    "    conn = connect()",
    "    cursor = conn.cursor()",
    "    cursor.execute(`"`"`"",
    "        SELECT * FROM students WHERE ",
    "        id LIKE ? OR name LIKE ? OR age LIKE ? OR room LIKE ?",
    "    `"`"`", (f'%{keyword}%', f'%{keyword}%', f'%{keyword}%', f'%{keyword}%'))",
    "    rows = cursor.fetchall()",
    "    conn.close()",
    "    return rows",
    "",
    "def add_room(room_no, capacity):",
    "    conn = connect()",
    "    cursor = conn.cursor()",
    "    cursor.execute(`"INSERT INTO rooms (room_no, capacity, occupants) VALUES (?, ?, 0)`", (room_no, capacity))",
    "    conn.commit()",
    "    conn.close()",
    "",
    "def delete_room(room_no):",
    "    conn = connect()",
    "    cursor = conn.cursor()",
    "    cursor.execute(`"DELETE FROM rooms WHERE room_no = ?`", (room_no,))",
    "    conn.commit()",
    "    conn.close()",
    "",
    "def view_rooms():",
    "    conn = connect()",
    "    cursor = conn.cursor()",
    "    cursor.execute(`"SELECT * FROM rooms`")",
    "    rows = cursor.fetchall()",
    "    conn.close()",
    "    return rows",
    "",
    "def search_room(keyword):",
    "    conn = connect()",
    "    cursor = conn.cursor()",
    "    cursor.execute(`"SELECT * FROM rooms WHERE room_no LIKE ?`", (f'%{keyword}%',))",
    "    rows = cursor.fetchall()",
    "    conn.close()",
    "    return rows",
    "",
    "def add_payment(student_id, amount, date):",
    "    conn = connect()",
    "    cursor = conn.cursor()",
    "    cursor.execute(`"INSERT INTO payments (student_id, amount, date) VALUES (?, ?, ?)`", (student_id, amount, date))",
    "    conn.commit()",
    "    conn.close()",
    "",
    "def delete_payment(payment_id):",
    "    conn = connect()",
    "    cursor = conn.cursor()",
    "    cursor.execute(`"DELETE FROM payments WHERE payment_id = ?`", (payment_id,))",
    "    conn.commit()",
    "    conn.close()",
    "",
    "def view_payments():",
    "    conn = connect()",
    "    cursor = conn.cursor()",
    "    cursor.execute(`"SELECT * FROM payments`")",
    "    rows = cursor.fetchall()",
    "    conn.close()",
    "    return rows",
    "",
    "def search_payment(keyword):",
    "    conn = connect()",
    "    cursor = conn.cursor()",
    "    cursor.execute(`"`"`"",
    "        SELECT * FROM payments WHERE ",
    "        student_id LIKE ? OR amount LIKE ? OR date LIKE ?",
    "    `"`"`", (f'%{keyword}%', f'%{keyword}%', f'%{keyword}%'))",
    "    rows = cursor.fetchall()",
    "    conn.close()",
    "    return rows"
)

$d = $word.ActiveDocument
$full = $lines -join "`r"
$r = $d.Range(0, $d.Content.End)
$r.Text = $full
